$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "alpha4F"

# Add new row 16, mirroring the pattern of row 15 (index/label pair) with 1s across C:M
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) # xlPasteFormats - copy the bold/border/centered style

$ws.Range("B16").Value = $ws.Range("B15").Value2

$ws.Range("C16:M16").Value = 1
